$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5, shifting existing rows 5-8 down to 6-9
$ws.Rows.Item(5).Insert()

# Fill in the new row with only the city name (missing data for the rest)
$ws.Range("A5").Value = "City Z"

# Update selection to match target state
$ws.Range("A6").Select()
